$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Guide")
$ws.Range("E2").Value = "x"
$ws.Range("E2").Borders.Item(7).LineStyle = 1
$ws.Range("E2").Borders.Item(7).Weight = -4138
$ws.Range("E2").Borders.Item(7).ColorIndex = 64
